$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "1.00", "589.74")
# that must stay plain text, matching the source inlineStr cells. A
# leading apostrophe forces Excel to treat the value as text (same as
# typing it in the UI) instead of silently coercing it to a number; the
# subsequent Style reset clears the "quote prefix" cell style that Excel
# applies as a side effect, keeping the cell formatting unchanged.

# Row 2
$ws.Range("D2").Value = '''70.374.45'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.43%  '

# Row 3
$ws.Range("D3").Value = '''3.609.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.06%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").Value = '''589.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.38%  '

# Row 6
$ws.Range("D6").Value = '''190.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.59%  '

# Row 7
$ws.Range("E7").Value = '  +1.87%  '

# Row 8
$ws.Range("D8").Value = '''3.597.60'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.96%  '

# Row 9
$ws.Range("E9").Value = '  -0.06%  '

# Row 10
$ws.Range("E10").Value = '  +0.62%  '

# Row 11
$ws.Range("E11").Value = '  +2.77%  '

# Row 12
$ws.Range("D12").Value = '''58.09'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.35%  '

# Row 13
$ws.Range("D13").Value = '''0.0000291'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.43%  '

# Row 14
$ws.Range("D14").Value = '''9.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.44%  '

# Row 15
$ws.Range("D15").Value = '''4.177.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.96%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '''3.611.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.07%  '

# Row 17
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").Value = '''19.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.84%  '

# Row 18
$ws.Range("D18").Value = '''70.274.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.40%  '

# Row 19
$ws.Range("D19").Value = '''12.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.73%  '

# Row 20
$ws.Range("E20").Value = '  +0.28%  '

# Row 21
$ws.Range("E21").Value = '  +4.22%  '

# Row 22
$ws.Range("D22").Value = '''495.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.21%  '

# Row 23
$ws.Range("D23").Value = '''17.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +16.93%  '

# Row 24
$ws.Range("D24").Value = '''5.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.19%  '

# Row 25
$ws.Range("D25").Value = '''4.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.45%  '

# Row 26
$ws.Range("D26").Value = '''90.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.40%  '

# Row 27
$ws.Range("D27").Value = '''3.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.36%  '

# Row 28
$ws.Range("D28").Value = '''11.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.58%  '

# Row 29
$ws.Range("D29").Value = '''9.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.93%  '

# Row 30
$ws.Range("D30").Value = '''32.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.67%  '

# Row 31
$ws.Range("D31").Value = '''7.60'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.00%  '

# Row 32
$ws.Range("D32").Value = '''12.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.28%  '

# Row 33
$ws.Range("D33").Value = '''617.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.95%  '

# Row 34
$ws.Range("D34").Value = '''0.118'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.97%  '

# Row 35
$ws.Range("D35").Value = '''65.23'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.00%  '

# Row 36
$ws.Range("D36").Value = '''0.0₃0819'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.91%  '

# Row 37
$ws.Range("B37").Value = 'TheGraph'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D37").Value = '''0.405'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.38%  '

# Row 38
$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D38").Value = '''38.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.18%  '

# Row 39
$ws.Range("E39").Value = '  +0.09%  '

# Row 40
$ws.Range("E40").Value = '  -1.02%  '

# Row 41
$ws.Range("D41").Value = '''3.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.38%  '

# Row 42
$ws.Range("D42").Value = '''3.307.92'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.77%  '

# Row 43
$ws.Range("D43").Value = '''3.09'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.70%  '

# Row 44
$ws.Range("E44").Value = '  +4.63%  '

# Row 45
$ws.Range("D45").Value = '''2.68'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.26%  '

# Row 46
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '''3.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.54%  '

# Row 47
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = '''0.138'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.25%  '

# Row 48
$ws.Range("D48").Value = '''9.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.71%  '

# Row 49
$ws.Range("D49").Value = '''2.71'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.70%  '

# Row 50
$ws.Range("D50").Value = '''3.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.18%  '

# Row 51
$ws.Range("D51").Value = '''0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '
